# Update gh-pages output numbers ("want to go" counts) for the latest scrape.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 8934
$ws1.Range("F7").Value  = 11485
$ws1.Range("F13").Value = 127
$ws1.Range("F28").Value = 1396
$ws1.Range("F35").Value = 477
$ws1.Range("F37").Value = 76

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value  = 354

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 354
$ws4.Range("F9").Value  = 8934
$ws4.Range("F11").Value = 11485
$ws4.Range("F28").Value = 1396
$ws4.Range("F35").Value = 477

$wb.Save()
